$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "section 0" question is being inserted as its own column, before the
# existing "preg_test_1" column (C). Columns("C").Insert() shifts every
# existing column at/after C one slot to the right (C->D, D->E, E->F),
# carrying cell values/types/styles and growing the used range/dimension
# along with it (A1:E7 -> A1:F7).
$ws.Columns("C").Insert()

# Header text for the freshly inserted column C.
$ws.Range("C1").Value = "preg_secc0"

# The header-row cell comments ("test 1" / "test crit" / "test 3") stay
# anchored to their original cell references instead of moving with the
# column insert, so re-home them by hand. Walk right-to-left so every
# source comment is read before its slot gets overwritten; the right-most
# slot (F1) never had a comment so it needs a brand-new one, and C1 gets the
# new section-0 text.
[void]$ws.Range("F1").AddComment($ws.Range("E1").Comment.Text())
[void]$ws.Range("E1").Comment.Text($ws.Range("D1").Comment.Text())
[void]$ws.Range("D1").Comment.Text($ws.Range("C1").Comment.Text())
[void]$ws.Range("C1").Comment.Text("testing secciones (seccion 0, debería ir primero)")
